$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-03 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-04 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("20÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "83÷8=", 2) | Out-Null
$d.Content.Find.Execute("11÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷5=", 2) | Out-Null
$d.Content.Find.Execute("96÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "28÷2=", 2) | Out-Null
$d.Content.Find.Execute("15÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷3=", 2) | Out-Null
$d.Content.Find.Execute("25÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷8=", 2) | Out-Null
$d.Content.Find.Execute("13÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷3=", 2) | Out-Null
$d.Content.Find.Execute("29÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "28÷9=", 2) | Out-Null
$d.Content.Find.Execute("48÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "60÷4=", 2) | Out-Null
$d.Content.Find.Execute("71÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "22÷8=", 2) | Out-Null
$d.Content.Find.Execute("26÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "96÷9=", 2) | Out-Null
$d.Content.Find.Execute("88÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "17÷7=", 2) | Out-Null
$d.Content.Find.Execute("84÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "11÷4=", 2) | Out-Null
$d.Content.Find.Execute("63÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷7=", 2) | Out-Null
$d.Content.Find.Execute("29÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷5=", 2) | Out-Null
$d.Content.Find.Execute("11÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷6=", 2) | Out-Null
$d.Content.Find.Execute("57÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "20÷8=", 2) | Out-Null
$d.Content.Find.Execute("63÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷3=", 2) | Out-Null
$d.Content.Find.Execute("34÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "89÷7=", 2) | Out-Null
$d.Content.Find.Execute("66÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "68÷4=", 2) | Out-Null
$d.Content.Find.Execute("87÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "69÷4=", 2) | Out-Null
$d.Content.Find.Execute("17÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷5=", 2) | Out-Null
$d.Content.Find.Execute("98÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "91÷6=", 2) | Out-Null
$d.Content.Find.Execute("91÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "34÷9=", 2) | Out-Null
$d.Content.Find.Execute("36÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷6=", 2) | Out-Null
$d.Content.Find.Execute("80÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "80÷7=", 2) | Out-Null
